$d = $word.ActiveDocument

# 1. Ultrasound machine model text change
$d.Content.Find.Execute("MEDISONIC MODELO H60 ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "MINDRAY MODELO DC – N3 ", 2)

# 2. Merge "Diámetro " + "Biparietal" + " " runs into a single run's text
$d.Content.Find.Execute("Diámetro Biparietal ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Diámetro Biparietal ", 2)

# 3. Merge " Lat. " + "x" + " min., registrado mediante " + "Doppler" + " pulsado y continuo en Modo "
$d.Content.Find.Execute(" Lat. x min., registrado mediante Doppler pulsado y continuo en Modo ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " Lat. x min., registrado mediante Doppler pulsado y continuo en Modo ", 2)

# 4. Merge "Muestra trayecto " + "espiralado" + " habitual. "
$d.Content.Find.Execute("Muestra trayecto espiralado habitual. ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Muestra trayecto espiralado habitual. ", 2)
